$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row for the 40-target section, reusing the existing
# "Target = 40" shared string (row 43, column A only).
$ws.Range("A43").Value = "Target = 40"

# New "20sec" calibration run at target 40 (row 44).
$ws.Range("A44").Value = "20sec"
$ws.Range("B44").Value = 437.6
$ws.Range("C44").Value = 434.4
$ws.Range("D44").Value = 439.8
$ws.Range("E44").Formula = "=(B44+C44+D44)/3"
$ws.Range("F44").Formula = "=E44/20"

# Slightly narrower first column after the new data was added.
$ws.Columns.Item(1).ColumnWidth = 21

# Scroll the view down and move the active selection to the new data,
# matching where the author left the cursor after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
$ws.Range("G44").Select() | Out-Null
